$d = $word.ActiveDocument

# --- 1. Swap the bracketed phone-confirmation placeholder for the real follow-up note ---
$old = "(IF PHONE CONFIRMATION: A follow up phone call to the site in the summer of 2020 confirmed that the development [WAS/WAS NOT] checking the site and removing waste seven days a week.1)"
$new = "The consolidation could not be reached for a follow up interview at this time to verify if implementation of AWS has allowed the to become compliant. "

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# --- 2. The paragraph now carries an explicit 8pt (160 twip) "space after" alongside its
#         existing auto-spacing, matching the edited copy's tighter layout. ---
$p3 = $d.Paragraphs.Item(3)
$p3.Format.SpaceAfter = 8

# --- 3. The blank paragraph right after it switches from the document-default spacing
#         scheme (after=160/line=259) to the same auto-spacing/single-line scheme used
#         by the paragraph above it. ---
$p4 = $d.Paragraphs.Item(4)
$p4.Format.SpaceBeforeAuto = 1
$p4.Format.SpaceAfterAuto = 1
$p4.Format.SpaceAfter = 0
$p4.Format.LineSpacingRule = 0
